$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the existing rows 6 and 7 down to 7 and 8.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly price entry.
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(6, 3).Value = "Los Lagos"
$ws.Cells.Item(6, 4).Value = 44523
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Value = 100112012
$ws.Cells.Item(6, 7).Value = "Espinaca"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 9000
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 9000
$ws.Cells.Item(6, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 900
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = "Hortaliza"
